$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.599562883377075
$ws.Range("B1").Value = 3.780060529708862
$ws.Range("C1").Value = 2.409914970397949
$ws.Range("D1").Value = 0.6090931296348572
$ws.Range("E1").Value = 0.9535717964172363
